$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated coin price/volume figures for the 11-2-2023 13:00 snapshot.
# Each cell must stay stored as TEXT (matches the source feed's inline-string
# cells), so we force the Text number format before writing the literal and
# reset the style afterwards to avoid leaving a stray style index behind.
$newValues = [ordered]@{
    "D2" = "309.21"
    "E2" = "0.87%"
    "D3" = "41.16"
    "E3" = "4.95%"
    "D4" = "5.116"
    "E4" = "1.02%"
    "D5" = "0.07631"
    "E5" = "-0.50%"
    "D6" = "1.622"
    "E6" = "0.86%"
    "E7" = "3.06%"
    "D8" = "0.9050"
    "E8" = "-1.11%"
    "D9" = "0.1109"
    "E9" = "7.88%"
    "D10" = "0.1795"
    "E10" = "3.05%"
    "D11" = "0.09047"
    "E11" = "1.26%"
    "D12" = "0.04255"
    "E12" = "-3.89%"
    "D13" = "0.1049"
    "E13" = "-0.48%"
    "E14" = "-0.48%"
    "D15" = "0.005671"
    "E15" = "-3.15%"
    "D16" = "3.341"
    "E16" = "-0.47%"
    "D17" = "4.248"
    "E17" = "0.42%"
    "E18" = "0.17%"
    "D19" = "6.674"
    "E19" = "-4.74%"
    "D20" = "0.1364"
    "E20" = "1.27%"
    "D21" = "0.2706"
    "E21" = "-1.14%"
    "D22" = "0.04029"
    "E22" = "-2.73%"
    "E23" = "4.32%"
    "D24" = "0.004093"
    "E24" = "0.15%"
    "E25" = "-0.12%"
    "D26" = "0.0003742"
    "D38" = "0.02410"
    "E38" = "1.73%"
    "D39" = "0.05241"
    "E39" = "1.56%"
    "D40" = "0.007757"
    "E40" = "-2.11%"
    "E41" = "-1.34%"
    "D42" = "0.007044"
    "E42" = "20.61%"
    "E43" = "-0.09%"
    "D44" = "0.008436"
    "E44" = "13.91%"
    "D45" = "0.3339"
    "E45" = "0.40%"
    "D46" = "0.00006893"
    "E46" = "7.20%"
    "E47" = "-0.15%"
    "D48" = "0.05937"
    "E48" = "1,323.07%"
    "E50" = "-0.15%"
    "E51" = "-0.15%"
}

foreach ($ref in $newValues.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $newValues[$ref]
    $cell.Style = "Normal"
}
